$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The Board summary table: "Stock" and "Waste" rows are merged into a
# single "Stock/Waste" row, and the old "Waste" row is reused to show
# "Foundation" (the Foundation row that used to follow it is removed).
$ws.Range("A38").Value = "Stock/Waste"
$ws.Range("A39").Value = "Foundation"
$ws.Range("B39").Value = "Foundation"
$ws.Range("A40").ClearContents()
$ws.Range("B40").ClearContents()

# The note about re-opening a tableau card moves out of the table's O
# column, gets reworded, and is placed as a new note below the table.
$ws.Range("O36").ClearContents()
$ws.Range("D41").Value = "If card was moved out of tableau pile a new one should be opened on the origin pile"

# New "CardCanBeMoved" field added to the Tableau class property list
# (row 28 was blank before, between IsClosed at row 27 and the Foundation
# header at row 29). It keeps the same green highlight style as row 27.
$ws.Range("A28").Value = "CardCanBeMoved"
$ws.Range("B28").Value = "bool"

# Reflect the updated scroll position and selection shown in the file.
$ws.Select()
$appWin = $excel.ActiveWindow
$appWin.ScrollRow = 14
$appWin.ScrollColumn = 1
$ws.Range("A27").Select() | Out-Null
